$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "uncertain"
$ws.Range("C3").Value = 0.08019999999999999

$ws.Range("B4").Value = "interest"
$ws.Range("C4").Value = -0.0263

$ws.Range("B5").Value = "inflation"
$ws.Range("C5").Value = -0.2366

$ws.Range("B6").Value = "trade"
$ws.Range("C6").Value = -0.0926

$ws.Range("B8").Value = "interest"
$ws.Range("C8").Value = -0.4211

$ws.Range("B9").Value = "uncertain"
$ws.Range("C9").Value = -0.6423

$ws.Range("B10").Value = "invest"
$ws.Range("C10").Value = -0.0548

$ws.Range("B11").Value = "trade"
$ws.Range("C11").Value = -0.469

$ws.Range("B12").Value = "uncertain"
$ws.Range("C12").Value = -0.1748

$ws.Range("B13").Value = "interest"
$ws.Range("C13").Value = 0.4966

$ws.Range("B14").Value = "trade"
$ws.Range("C14").Value = 0.0533

$ws.Range("B17").Value = "inflation"
$ws.Range("C17").Value = 0.5832000000000001

$ws.Range("B18").Value = "interest"
$ws.Range("C18").Value = 0.8001

$ws.Range("B19").Value = "invest"
$ws.Range("C19").Value = -0.1405

$ws.Range("B20").Value = "trade"
$ws.Range("C20").Value = 0.6494
